# Updated cryptos list with GitHub Actions scrape refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume(1h) (E) columns hold plain text in this sheet (e.g. "3.080",
# "1.000"), so force the text format before writing new values -- otherwise Excel
# would auto-convert numeric-looking strings to numbers and drop trailing zeros.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "23.764.04"
$ws.Range("E2").Value = "  +1.96%  "
$ws.Range("D3").Value = "1.654.46"
$ws.Range("E3").Value = "  +1.87%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "1.001"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("D6").Value = "304.12"
$ws.Range("E6").Value = "  +0.47%  "
$ws.Range("D7").Value = "0.3814"
$ws.Range("E7").Value = "  +1.59%  "
$ws.Range("D8").Value = "0.3613"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").Value = "51.11"
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("E10").Value = "  +2.80%  "
$ws.Range("D11").Value = "0.08203"
$ws.Range("E11").Value = "  +0.89%  "
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").Value = "22.69"
$ws.Range("E13").Value = "  +2.05%  "
$ws.Range("D14").Value = "6.538"
$ws.Range("E14").Value = "  +1.25%  "
$ws.Range("D15").Value = "7.431"
$ws.Range("E15").Value = "  +2.28%  "
$ws.Range("D16").Value = "0.00001239"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("D17").Value = "1.637.75"
$ws.Range("E17").Value = "  +1.40%  "
$ws.Range("D18").Value = "97.84"
$ws.Range("E18").Value = "  +4.04%  "
$ws.Range("D19").Value = "0.06969"
$ws.Range("E19").Value = "  +0.51%  "
$ws.Range("D20").Value = "6.783"
$ws.Range("E20").Value = "  +3.53%  "
$ws.Range("D21").Value = "17.74"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "12.72"
$ws.Range("E23").Value = "  +2.16%  "
$ws.Range("D24").Value = "23.769.07"
$ws.Range("E24").Value = "  +1.96%  "
$ws.Range("D25").Value = "2.553"
$ws.Range("E25").Value = "  +2.40%  "
$ws.Range("D26").Value = "3.080"
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "21.32"
$ws.Range("E27").Value = "  +1.07%  "
$ws.Range("D28").Value = "151.07"
$ws.Range("E28").Value = "  +0.41%  "
$ws.Range("D29").Value = "5.215"
$ws.Range("E29").Value = "  -1.37%  "
$ws.Range("D30").Value = "134.98"
$ws.Range("E30").Value = "  +1.65%  "
$ws.Range("D31").Value = "1.826.47"
$ws.Range("E31").Value = "  +1.69%  "
$ws.Range("D32").Value = "6.914"
$ws.Range("E32").Value = "  +3.07%  "
$ws.Range("D33").Value = "1.080"
$ws.Range("E33").Value = "  +2.05%  "
$ws.Range("D34").Value = "2.144"
$ws.Range("E34").Value = "  +2.82%  "
$ws.Range("D35").Value = "11.94"
$ws.Range("E35").Value = "  +6.52%  "
$ws.Range("D36").Value = "0.02833"
$ws.Range("E36").Value = "  +3.43%  "
$ws.Range("D37").Value = "0.2518"
$ws.Range("E37").Value = "  +1.84%  "
$ws.Range("D38").Value = "6.141"
$ws.Range("E38").Value = "  +2.83%  "
$ws.Range("D39").Value = "0.08834"
$ws.Range("E39").Value = "  +0.84%  "
$ws.Range("D40").Value = "0.07179"
$ws.Range("E40").Value = "  +1.32%  "
$ws.Range("D41").Value = "13.13"
$ws.Range("E41").Value = "  +9.42%  "
$ws.Range("D42").Value = "0.7068"
$ws.Range("E42").Value = "  +1.67%  "
$ws.Range("D43").Value = "1.342"
$ws.Range("E43").Value = "  +1.16%  "
$ws.Range("D44").Value = "15.90"
$ws.Range("E44").Value = "  +0.58%  "
$ws.Range("D45").Value = "0.6542"
$ws.Range("E45").Value = "  +1.69%  "
$ws.Range("D46").Value = "2.330"
$ws.Range("E46").Value = "  +3.02%  "
$ws.Range("D47").Value = "1.000"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("D48").Value = "3.964"
$ws.Range("E48").Value = "  +0.27%  "
$ws.Range("D49").Value = "0.07982"
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("D50").Value = "128.56"
$ws.Range("E50").Value = "  +2.23%  "
$ws.Range("D51").Value = "1.194"
$ws.Range("E51").Value = "  +1.00%  "
